$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D3").Value = "2016-01-19 06:54:44"
$wsZhCn.Range("G3").Value = "2016-01-19 06:55:28"

$wsDeDe.Range("D3").Value = "2016-01-19 06:54:54"
$wsDeDe.Range("G3").Value = "2016-01-19 06:55:46"
